$wb = $excel.ActiveWorkbook

# ---- Participants ----
$ws = $wb.Worksheets.Item("Participants")
$ws.Columns("D:D").Insert()
$ws.Range("D1").Value = "Rank-ICPC 2025 Team Formation - 03.xlsx"

$ws.Range("A2").Value = "Tamjid_Hossen(Tamjid)"
$ws.Range("B2").Value = 200
$ws.Range("C2").Value = 300
$ws.Range("D2").Value = 300
$ws.Range("E2").Value = 800
$ws.Range("A3").Value = "YouDOntKNowWHo(Nabeel Ahsan)"
$ws.Range("B3").Value = 225
$ws.Range("C3").Value = 258
$ws.Range("D3").Value = 225
$ws.Range("E3").Value = 708
$ws.Range("A4").Value = "shazidmashrafi(Shazid)"
$ws.Range("B4").Value = 180
$ws.Range("C4").Value = 225
$ws.Range("D4").Value = 258
$ws.Range("E4").Value = 663
$ws.Range("A5").Value = "sf61561(Syed Fahad Mahmud)"
$ws.Range("B5").Value = 300
$ws.Range("C5").Value = 164
$ws.Range("D5").Value = 164
$ws.Range("E5").Value = 628
$ws.Range("A6").Value = "y_this_kolaveri(SAIF)"
$ws.Range("B6").Value = 258
$ws.Range("C6").Value = 120
$ws.Range("D6").Value = 129
$ws.Range("E6").Value = 507
$ws.Range("A7").Value = "Md_Saurob_bhuyan(Noob)"
$ws.Range("B7").Value = 164
$ws.Range("C7").Value = 139
$ws.Range("D7").Value = 180
$ws.Range("E7").Value = 483
$ws.Range("A8").Value = "rakin_ahsan(Rakin)"
$ws.Range("B8").Value = 106
$ws.Range("C8").Value = 180
$ws.Range("D8").Value = 120
$ws.Range("E8").Value = 406
$ws.Range("A9").Value = "farhanshadiq(Farhan)"
$ws.Range("B9").Value = 113
$ws.Range("C9").Value = 129
$ws.Range("D9").Value = 139
$ws.Range("E9").Value = 381
$ws.Range("A10").Value = "Aniksamiul(Anik)"
$ws.Range("B10").Value = 120
$ws.Range("C10").Value = 106
$ws.Range("D10").Value = 150
$ws.Range("E10").Value = 376
$ws.Range("A11").Value = "Akash_khan"
$ws.Range("B11").Value = 139
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 200
$ws.Range("E11").Value = 339
$ws.Range("A12").Value = "AL_AMIN_17(Al Amin)"
$ws.Range("B12").Value = 129
$ws.Range("C12").Value = 200
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 329
$ws.Range("A13").Value = "Apon_Chowdhury(Apon)"
$ws.Range("B13").Value = 100
$ws.Range("C13").Value = 113
$ws.Range("D13").Value = 113
$ws.Range("E13").Value = 326
$ws.Range("A14").Value = "Noornabi1770(Noor)"
$ws.Range("B14").Value = 90
$ws.Range("C14").Value = 150
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 240
$ws.Range("A15").Value = "Marufhussain(maruf)"
$ws.Range("B15").Value = 150
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 150
$ws.Range("A16").Value = "SadmanIshtiak(Sadman)"
$ws.Range("B16").Value = 95
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 95
$ws.Range("A17").Value = "_Mohiul007(Rabby)"
$ws.Range("B17").Value = 86
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 86

# ---- Team_1 ----
$ws = $wb.Worksheets.Item("Team_1")
$ws.Columns("D:D").Insert()
$ws.Range("D1").Value = "Rank-ICPC 2025 Team Formation - 03.xlsx"

$ws.Range("A2").Value = "Tamjid_Hossen(Tamjid)"
$ws.Range("B2").Value = 200
$ws.Range("C2").Value = 300
$ws.Range("D2").Value = 300
$ws.Range("E2").Value = 800
$ws.Range("A3").Value = "YouDOntKNowWHo(Nabeel Ahsan)"
$ws.Range("B3").Value = 225
$ws.Range("C3").Value = 258
$ws.Range("D3").Value = 225
$ws.Range("E3").Value = 708
$ws.Range("A4").Value = "shazidmashrafi(Shazid)"
$ws.Range("B4").Value = 180
$ws.Range("C4").Value = 225
$ws.Range("D4").Value = 258
$ws.Range("E4").Value = 663

# ---- Team_2 ----
$ws = $wb.Worksheets.Item("Team_2")
$ws.Columns("D:D").Insert()
$ws.Range("D1").Value = "Rank-ICPC 2025 Team Formation - 03.xlsx"

$ws.Range("A2").Value = "sf61561(Syed Fahad Mahmud)"
$ws.Range("B2").Value = 300
$ws.Range("C2").Value = 164
$ws.Range("D2").Value = 164
$ws.Range("E2").Value = 628
$ws.Range("A3").Value = "y_this_kolaveri(SAIF)"
$ws.Range("B3").Value = 258
$ws.Range("C3").Value = 120
$ws.Range("D3").Value = 129
$ws.Range("E3").Value = 507
$ws.Range("A4").Value = "Md_Saurob_bhuyan(Noob)"
$ws.Range("B4").Value = 164
$ws.Range("C4").Value = 139
$ws.Range("D4").Value = 180
$ws.Range("E4").Value = 483

# ---- Team_3 ----
$ws = $wb.Worksheets.Item("Team_3")
$ws.Columns("D:D").Insert()
$ws.Range("D1").Value = "Rank-ICPC 2025 Team Formation - 03.xlsx"

$ws.Range("A2").Value = "rakin_ahsan(Rakin)"
$ws.Range("B2").Value = 106
$ws.Range("C2").Value = 180
$ws.Range("D2").Value = 120
$ws.Range("E2").Value = 406
$ws.Range("A3").Value = "farhanshadiq(Farhan)"
$ws.Range("B3").Value = 113
$ws.Range("C3").Value = 129
$ws.Range("D3").Value = 139
$ws.Range("E3").Value = 381
$ws.Range("A4").Value = "Aniksamiul(Anik)"
$ws.Range("B4").Value = 120
$ws.Range("C4").Value = 106
$ws.Range("D4").Value = 150
$ws.Range("E4").Value = 376

# ---- Team_4 ----
$ws = $wb.Worksheets.Item("Team_4")
$ws.Columns("D:D").Insert()
$ws.Range("D1").Value = "Rank-ICPC 2025 Team Formation - 03.xlsx"

$ws.Range("A2").Value = "Akash_khan"
$ws.Range("B2").Value = 139
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 200
$ws.Range("E2").Value = 339
$ws.Range("A3").Value = "AL_AMIN_17(Al Amin)"
$ws.Range("B3").Value = 129
$ws.Range("C3").Value = 200
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 329
$ws.Range("A4").Value = "Apon_Chowdhury(Apon)"
$ws.Range("B4").Value = 100
$ws.Range("C4").Value = 113
$ws.Range("D4").Value = 113
$ws.Range("E4").Value = 326

# ---- Team_5 ----
$ws = $wb.Worksheets.Item("Team_5")
$ws.Columns("D:D").Insert()
$ws.Range("D1").Value = "Rank-ICPC 2025 Team Formation - 03.xlsx"

$ws.Range("A2").Value = "Noornabi1770(Noor)"
$ws.Range("B2").Value = 90
$ws.Range("C2").Value = 150
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 240
$ws.Range("A3").Value = "Marufhussain(maruf)"
$ws.Range("B3").Value = 150
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 150
$ws.Range("A4").Value = "SadmanIshtiak(Sadman)"
$ws.Range("B4").Value = 95
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 95

# ---- Team_6 ----
$ws = $wb.Worksheets.Item("Team_6")
$ws.Columns("D:D").Insert()
$ws.Range("D1").Value = "Rank-ICPC 2025 Team Formation - 03.xlsx"

$ws.Range("A2").Value = "_Mohiul007(Rabby)"
$ws.Range("B2").Value = 86
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 86
